$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1. Title line: "Numerical An" + (bookmark _GoBack) + "alysis of
#    Diﬀerential Equations using " -> single run
#    "Numerical Analysis of Diﬀerential Equations using "
#    A same-text Find/Replace that spans the old run break (and the
#    bookmark sitting on that break) merges the two runs into one and
#    drops the bookmark that was anchored there, exactly as typing
#    over/through it in Word would.
# -----------------------------------------------------------------
$titleOld = "Numerical An" + "alysis of Di" + [char]0xFB00 + "erential Equations using "
$titleNew = "Numerical An" + "alysis of Di" + [char]0xFB00 + "erential Equations using "
$d.Content.Find.Execute($titleOld, $true, $false, $false, $false, $false, $true, 1, $false, $titleNew, 2) | Out-Null

# -----------------------------------------------------------------
# 2. Author credit line: " rg517, Alex Constantinou " + "acX17"(yellow)
#    -> " rg517, Alex Constantinou" + " ac517" (no highlight)
# -----------------------------------------------------------------

# 2a. Locate the run containing "acX17" and rewrite its text to include
#     a leading space and the corrected username, keeping it a
#     separate run from the text before it.
$acRange = $d.Content
$acRange.Find.Execute("acX17", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$acStart = $acRange.Start
$acEnd = $acRange.End
$acRange2 = $d.Range($acStart, $acEnd)
$acRange2.Text = " ac517"

# 2b. Remove the now-redundant trailing space that used to separate
#     "Constantinou" from "acX17" (it sat immediately before $acStart).
$spaceRange = $d.Range($acStart - 1, $acStart)
$spaceRange.Delete()

# 2c. Strip the yellow highlight that used to mark "acX17" from the
#     (now relocated) "ac517" run. Use a formatting-only find/replace
#     scoped to a Range built straight from the offsets recorded above
#     (no intervening text search) so the Find engine's format-replace
#     state is not disturbed.
$hiRange = $d.Range($acStart - 1, $acEnd)
$f = $hiRange.Find
$f.ClearFormatting()
$f.Replacement.ClearFormatting()
$f.Text = ""
$f.Replacement.Text = ""
$f.Highlight = $true
$f.Replacement.Highlight = $false
$f.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null

# -----------------------------------------------------------------
# 3. Re-anchor the "_GoBack" bookmark (which the title edit above
#    removed) on the empty paragraph right before the final sectPr -
#    i.e. where the cursor was left after finishing Alex's username.
# -----------------------------------------------------------------
$lastPara = $d.Paragraphs.Last.Range
$d.Bookmarks.Add("_GoBack", $lastPara)
